# Update the "Price" column (D) with refreshed symbol-list quotes.
# Values are written with a leading apostrophe so Excel keeps them as
# text (matching the workbook's existing inline-string cell type)
# instead of auto-converting the numeric-looking text into a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.89"
$ws.Range("D3").Value = "'23.07"
$ws.Range("D4").Value = "'5.416"
$ws.Range("D5").Value = "'0.05943"
$ws.Range("D6").Value = "'3.392"
$ws.Range("D7").Value = "'0.8068"
$ws.Range("D9").Value = "'0.1414"
$ws.Range("D10").Value = "'0.07396"
$ws.Range("D12").Value = "'0.03079"
$ws.Range("D13").Value = "'0.09347"
$ws.Range("D14").Value = "'3.933"
$ws.Range("D15").Value = "'0.001587"
$ws.Range("D18").Value = "'0.005465"
$ws.Range("D20").Value = "'0.0009820"
$ws.Range("D21").Value = "'0.00007505"
$ws.Range("D22").Value = "'3.655"
$ws.Range("D23").Value = "'6.446"
$ws.Range("D25").Value = "'0.3243"
$ws.Range("D26").Value = "'0.1341"
$ws.Range("D40").Value = "'0.03911"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D44").Value = "'0.007296"
$ws.Range("D45").Value = "'0.00005174"
$ws.Range("D46").Value = "'0.00000000750"
